$d = $word.ActiveDocument

$replacements = @(
    @{old="34×48="; new="36×30="},
    @{old="61×24="; new="21×47="},
    @{old="52×29="; new="14×15="},
    @{old="91×35="; new="81×53="},
    @{old="60×16="; new="49×38="},
    @{old="89×17="; new="99×60="},
    @{old="51×72="; new="88×12="},
    @{old="19×82="; new="28×57="},
    @{old="36×16="; new="39×88="},
    @{old="20×54="; new="37×53="},
    @{old="62×81="; new="71×96="},
    @{old="65×79="; new="69×96="},
    @{old="44×12="; new="63×22="},
    @{old="87×55="; new="71×93="},
    @{old="91×71="; new="89×36="},
    @{old="29×82="; new="39×76="},
    @{old="12×88="; new="70×70="},
    @{old="79×87="; new="18×82="},
    @{old="63×18="; new="36×86="},
    @{old="84×81="; new="83×78="},
    @{old="14×18="; new="49×54="},
    @{old="87×57="; new="12×80="},
    @{old="85×42="; new="76×31="},
    @{old="74×88="; new="66×22="},
    @{old="41×46="; new="80×67="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
